$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-09-19 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-09-20 Friday", 2)

# Update the division problems in the table, addressed by row/column
# to avoid ambiguity since several old/new values repeat.
$t = $d.Tables.Item(1)
$t.Cell(1, 1).Range.Text = "75÷8="
$t.Cell(1, 2).Range.Text = "53÷8="
$t.Cell(1, 3).Range.Text = "81÷4="
$t.Cell(1, 4).Range.Text = "82÷5="
$t.Cell(1, 5).Range.Text = "27÷6="
$t.Cell(5, 1).Range.Text = "95÷2="
$t.Cell(5, 2).Range.Text = "77÷5="
$t.Cell(5, 3).Range.Text = "48÷9="
$t.Cell(5, 4).Range.Text = "42÷8="
$t.Cell(5, 5).Range.Text = "45÷2="
$t.Cell(9, 1).Range.Text = "24÷3="
$t.Cell(9, 2).Range.Text = "27÷7="
$t.Cell(9, 3).Range.Text = "39÷2="
$t.Cell(9, 4).Range.Text = "82÷7="
$t.Cell(9, 5).Range.Text = "42÷2="
$t.Cell(13, 1).Range.Text = "69÷4="
$t.Cell(13, 2).Range.Text = "67÷3="
$t.Cell(13, 3).Range.Text = "65÷2="
$t.Cell(13, 4).Range.Text = "41÷9="
$t.Cell(13, 5).Range.Text = "51÷7="
$t.Cell(17, 1).Range.Text = "23÷6="
$t.Cell(17, 2).Range.Text = "14÷7="
$t.Cell(17, 3).Range.Text = "43÷3="
$t.Cell(17, 4).Range.Text = "82÷9="
$t.Cell(17, 5).Range.Text = "77÷7="
